# Insert a new weekly price-report row for "Espinaca" (row 266) on the
# "Terminal La Palmera de La Serena" sheet. All rows from the old 266
# onward shift down by one (handled automatically by Rows.Insert), and
# we populate the freshly inserted row with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 266, pushing existing row 266..398 down to 267..399
$ws.Rows.Item(266).Insert()

$newRow = 266

$ws.Cells.Item($newRow, 1).Value = 8
$ws.Cells.Item($newRow, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($newRow, 3).Value = "Coquimbo"

$ws.Cells.Item($newRow, 4).Value = 45016
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($newRow + 1, 4).NumberFormat

$ws.Cells.Item($newRow, 5).Value = 4
$ws.Cells.Item($newRow, 6).Value = 100112012
$ws.Cells.Item($newRow, 7).Value = "Espinaca"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 1600
$ws.Cells.Item($newRow, 11).Value = 400
$ws.Cells.Item($newRow, 12).Value = 500
$ws.Cells.Item($newRow, 13).Value = 450
$ws.Cells.Item($newRow, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item($newRow, 15).Value = "Provincia del Elquí"
$ws.Cells.Item($newRow, 16).Value = 900
$ws.Cells.Item($newRow, 17).Value = 0.5
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
